$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "Status" column (H) had every "won" result renamed to "paid".
$xlWhole = -4163
$statusRange = $ws.Range("H2:H21")
$statusRange.Replace("won", "paid", $xlWhole, 1, $false, $false, $true)

# The Status column (H) now shows an explicit best-fit width after the edit.
$ws.Columns.Item(8).AutoFit()

# The user's final selection landed on cell J11.
$ws.Range("J11").Select()
